$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new (blank) row above row 13. This shifts the old rows 13..23
#    down to 14..24, carrying their formatting/row-heights with them.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# The newly inserted row 13 picked up column-A formatting from the row above
# (Insert() clones the row above's look). The target layout has NO cell in
# A13 at all, so clear it out completely.
$ws.Range("A13").Clear()

# Populate the new row 13 (B13/C13) with the "Docentes responsaveis" value
# that used to live in the old row 10 (B10/C10), matching the column B/C
# formatting used throughout the sheet (copy format from B10/C10, which are
# still untouched at this point).
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "5840535 - Messias Borges Silva"
$ws.Range("C13").Value = "5840535 - Messias Borges Silva"

# ---------------------------------------------------------------------------
# 2) Update the cell values (text only -- formatting/heights are already
#    correct after the row-insert shift above).
# ---------------------------------------------------------------------------

# Row 10 "Objetivos:" -- new objectives paragraph.
$ws.Range("B10").Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre gestão da qualidade."
$ws.Range("C10").Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre gestão da qualidade."

# Row 14 (old row 13) "Programa resumido:" -- new short-syllabus text.
$ws.Range("B14").Value = "A definir, de acordo com o tópico programado."
$ws.Range("C14").Value = "A definir, de acordo com o tópico programado."

# Row 16 (old row 15) "Programa:" -- new syllabus text.
$ws.Range("B16").Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia."
$ws.Range("C16").Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia."

# Row 19 (old row 18) "Método:" -- text that used to sit under "Critério:".
$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# Row 20 (old row 19) "Critério:" -- text that used to sit under "Norma de recuperação:".
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."

# Row 21 (old row 20) "Norma de recuperação:" -- text that used to sit under "Bibliografia:".
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# Row 22 (old row 21) "Bibliografia:" -- brand-new bibliography text.
$ws.Range("B22").Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas na área de gestão e produção."
$ws.Range("C22").Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas na área de gestão e produção."
